$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: re-apply the existing named styles by setting the underlying
# formatting properties (matches the workbook's existing style records so
# Excel reuses the same style index instead of allocating a new one).
# ---------------------------------------------------------------------------

function Set-MTitleStyle($rng) {
    # mtitleStyle: bold black "Century" 12, thin box border, centered
    $rng.Font.Name = "Century"
    $rng.Font.Size = 12
    $rng.Font.Bold = $true
    $rng.Font.Color = 0
    $rng.HorizontalAlignment = -4108
    $rng.Borders.LineStyle = 1
}

function Set-NormalStyle($rng) {
    # normalStyle: green "Century" 12, thin box border, centered
    $rng.Font.Name = "Century"
    $rng.Font.Size = 12
    $rng.Font.Bold = $false
    $rng.Font.Color = 32768
    $rng.HorizontalAlignment = -4108
    $rng.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): recomputed marks after handling float input.
# ---------------------------------------------------------------------------

Set-MTitleStyle $ws.Range("A10")
Set-MTitleStyle $ws.Range("A11")
Set-MTitleStyle $ws.Range("A12")

$ws.Range("B10").Value = 13
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 52
$ws.Range("E12").Value = "52/112"

# ---------------------------------------------------------------------------
# Drop the third Student Ans / Correct Ans block (columns G:H) entirely.
# ---------------------------------------------------------------------------

$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------------
# Second Student Ans / Correct Ans block (columns D:E) only keeps rows
# 16-18; the rest of the block is dropped.
# ---------------------------------------------------------------------------

$ws.Range("D16").Value = $ws.Range("E16").Value2
Set-NormalStyle $ws.Range("D16")

$ws.Range("D17").Value = $ws.Range("E17").Value2
Set-NormalStyle $ws.Range("D17")

$ws.Range("D18").Value = $ws.Range("E18").Value2
Set-NormalStyle $ws.Range("D18")

$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# First Student Ans / Correct Ans block (columns A:B): fill in the student
# answer whenever it matches the correct answer.
# ---------------------------------------------------------------------------

$matchingRows = 17, 19, 20, 22, 25, 27, 31, 32, 39, 40
foreach ($r in $matchingRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $ws.Cells.Item($r, 2).Value2
    Set-NormalStyle $cell
}
